# Adds "green hydrogen" and "low carbon hydrogen" as two new fuel columns
# (L and M) on the HPtFM sheet, mirroring the existing fuel columns:
#  - header cells in row 1 get the new fuel names
#  - rows 2-6 (pathway rows) get 0 for the new fuels
#  - rows 7-8 (derived pathway rows) feed the new columns through from the
#    rows they reference, same as the existing columns do

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HPtFM")

# New header labels in L1 / M1
$ws.Range("L1").Value() = "green hydrogen"
$ws.Range("M1").Value() = "low carbon hydrogen"

# Match the header formatting (right-aligned, wrap text) used by the rest
# of row 1's header cells
$ws.Range("K1").Copy()
$ws.Range("L1:M1").PasteSpecial(-4122)

# New data values (rows 2-6) default to 0, matching the other fuel columns
$ws.Range("L2:M6").Value() = 0

# Rows 7 and 8 pull the new columns through from the rows they mirror,
# exactly like the existing shared formulas in columns C:K
$ws.Range("L7:M7").Formula() = "=L2"
$ws.Range("L8:M8").Formula() = "=L3"

# Keep the new columns the same width as the rest of the data columns
$ws.Columns.Item(12).ColumnWidth() = $ws.Columns.Item(10).ColumnWidth()
$ws.Columns.Item(13).ColumnWidth() = $ws.Columns.Item(10).ColumnWidth()

# Match the selection left behind in the authored workbook
$ws.Range("L2").Select() | Out-Null
